# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 31
$ws.Range("F3").Value = 1200
$ws.Range("G3").Value = "不可售"
$ws.Range("F7").Value = 953
$ws.Range("F10").Value = 549
$ws.Range("F11").Value = 1422
$ws.Range("F13").Value = 1323
$ws.Range("F14").Value = 2980
$ws.Range("F15").Value = 383
$ws.Range("F16").Value = 1600
$ws.Range("F17").Value = 1350
$ws.Range("F18").Value = 780
$ws.Range("F20").Value = 1353
$ws.Range("F21").Value = 258
$ws.Range("F23").Value = 1112
$ws.Range("F25").Value = 3432
$ws.Range("F27").Value = 558
$ws.Range("F28").Value = 1522

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 47
$ws.Range("F8").Value = 20
$ws.Range("F12").Value = 73

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 792

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 31
$ws.Range("F3").Value = 792
$ws.Range("F4").Value = 1200
$ws.Range("G4").Value = "不可售"
$ws.Range("F12").Value = 7
$ws.Range("F13").Value = 47
$ws.Range("F15").Value = 20
$ws.Range("F17").Value = 953
$ws.Range("F20").Value = 549
$ws.Range("F21").Value = 1422
$ws.Range("F23").Value = 1323
$ws.Range("F24").Value = 2980
$ws.Range("F25").Value = 383
$ws.Range("F26").Value = 1600
$ws.Range("F27").Value = 1350
$ws.Range("F28").Value = 780
$ws.Range("F30").Value = 1353
$ws.Range("F31").Value = 258
$ws.Range("F35").Value = 1112
$ws.Range("F37").Value = 3432
$ws.Range("F39").Value = 558
$ws.Range("F40").Value = 1522
$ws.Range("F41").Value = 73
